# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update total "Valor Mora" summary (E11)
$ws.Range("E11").Value = 715112

# Data table rows 16-26: Tipo Doc (B) stays "CC"; N Doc Trabajador (C), Nombre
# Trabajador (D), Periodo Mora (E), Valor Mora (F), Salario Basico (G)
$rows = @(
    @{ Row = 16; Doc = "9146138";     Nombre = "EPARQUIO LUIS AMARANTO BELTRAN";      Mora = 180000; Salario = 4500000 },
    @{ Row = 17; Doc = "73266013";    Nombre = "WILLIAM ORTIZ ORTIZ";                  Mora = 56940;  Salario = 1423500 },
    @{ Row = 18; Doc = "19249925";    Nombre = "JAIME LUIS AMARANTO MERCADO";           Mora = 56940;  Salario = 1423500 },
    @{ Row = 19; Doc = "73265981";    Nombre = "RAMON OLIVERO SARA";                    Mora = 56940;  Salario = 1423500 },
    @{ Row = 20; Doc = "1045309902";  Nombre = "KEIDER DE JESUS SARA CASSIANI";          Mora = 32266;  Salario = 1423500 },
    @{ Row = 21; Doc = "72250907";    Nombre = "YESID ENRIQUE DE LA CRUZ ESCORCIA";      Mora = 72000;  Salario = 1800000 },
    @{ Row = 22; Doc = "1051358115";  Nombre = "JAVIER ENRIQUE SALAS ESCORCIA";          Mora = 56940;  Salario = 1423500 },
    @{ Row = 23; Doc = "73267118";    Nombre = "LUIS ALBERTO SARA CASSIANI";             Mora = 32266;  Salario = 1423500 },
    @{ Row = 24; Doc = "1051361081";  Nombre = "LEONAL JOSE CARRASQUILLA SARMIENTO";     Mora = 56940;  Salario = 1423500 },
    @{ Row = 25; Doc = "1127948009";  Nombre = "MANUEL JOSE REALES ORTIZ";               Mora = 56940;  Salario = 1423500 },
    @{ Row = 26; Doc = "7912541";     Nombre = "HENIO JOSE REALES TEJEDA";               Mora = 56940;  Salario = 1423500 }
)

foreach ($r in $rows) {
    $ws.Range("C" + $r.Row).Value = $r.Doc
    $ws.Range("D" + $r.Row).Value = $r.Nombre
    $ws.Range("E" + $r.Row).Value = "2507"
    $ws.Range("F" + $r.Row).Value = $r.Mora
    $ws.Range("G" + $r.Row).Value = $r.Salario
}
